$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: Volume/Number and report week dates (new weekly collection)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/18/2024  Through  11/24/2024"

# ---------------------------------------------------------------------------
# Helper pattern used throughout: when a cell needs to change its underlying
# type (numeric <-> text placeholder), first assign the new value, then copy
# number-format/style from a cell that already carries the desired style so
# the pasted-in cell matches the target style index, then (re)apply the
# value so the stored type (text vs number) is correct.
# -4122 == xlPasteFormats
# ---------------------------------------------------------------------------

# Row 15 - Rape
$ws.Range("F15").Value = 8
$ws.Range("I15").Value = 25
$ws.Range("K15").Value = 66.666666666666
$ws.Range("L15").Value = 4.166666666666
$ws.Range("M15").Value = 31.578947368421
$ws.Range("N15").Value = -30.555555555555

# Row 16 - Robbery (C16/D16 become "0" placeholders, E16 becomes "***.*")
$ws.Range("C16").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("D16").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("D16").PasteSpecial(-4122)

$ws.Range("E16").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E16").PasteSpecial(-4122)

$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -42.857142857142
$ws.Range("L16").Value = -3.428571428571
$ws.Range("M16").Value = -23.873873873873
$ws.Range("N16").Value = -82.172995780590

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = -20
$ws.Range("I17").Value = 294
$ws.Range("J17").Value = 257
$ws.Range("K17").Value = 14.396887159533
$ws.Range("L17").Value = 5.755395683453
$ws.Range("M17").Value = 68.965517241379
$ws.Range("N17").Value = -37.974683544303

# Row 18 - Burglary
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 10
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 152
$ws.Range("J18").Value = 119
$ws.Range("K18").Value = 27.731092436974
$ws.Range("L18").Value = -34.199134199134
$ws.Range("M18").Value = -33.333333333333
$ws.Range("N18").Value = -85.537583254043

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -81.25
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = -58.928571428571
$ws.Range("I19").Value = 423
$ws.Range("J19").Value = 519
$ws.Range("K19").Value = -18.497109826589
$ws.Range("L19").Value = -42.213114754098
$ws.Range("M19").Value = 24.778761061946
$ws.Range("N19").Value = 0.475059382422

# Row 20 - G.L.A.
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 139
$ws.Range("J20").Value = 164
$ws.Range("K20").Value = -15.243902439024
$ws.Range("L20").Value = 6.106870229007
$ws.Range("M20").Value = 27.522935779816
$ws.Range("N20").Value = -84.364454443194

# Row 21 - TOTAL
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -52
$ws.Range("F21").Value = 77
$ws.Range("G21").Value = 112
$ws.Range("H21").Value = -31.25
$ws.Range("I21").Value = 1205
$ws.Range("J21").Value = 1218
$ws.Range("K21").Value = -1.067323481116
$ws.Range("L21").Value = -23.443456162642
$ws.Range("M21").Value = 9.845031905195
$ws.Range("N21").Value = -68.587069864442

# Row 22 - Transit (C22/D22/E22 switch from placeholders to real numbers)
$ws.Range("C17").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 2

$ws.Range("D17").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1

$ws.Range("E17").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = 100

$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 22
$ws.Range("J22").Value = 33
$ws.Range("K22").Value = -33.333333333333
$ws.Range("L22").Value = -31.25
$ws.Range("M22").Value = -15.384615384615

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 65
$ws.Range("G24").Value = 71
$ws.Range("H24").Value = -8.450704225352
$ws.Range("I24").Value = 866
$ws.Range("J24").Value = 1117
$ws.Range("K24").Value = -22.470904207699
$ws.Range("L24").Value = -29.363784665579
$ws.Range("M24").Value = 20.781032078103

# Row 25 - Retail Theft
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -16.666666666666
$ws.Range("F25").Value = 32
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = -21.951219512195
$ws.Range("I25").Value = 355
$ws.Range("J25").Value = 634
$ws.Range("K25").Value = -44.006309148265
$ws.Range("L25").Value = -47.407407407407

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 120
$ws.Range("F26").Value = 54
$ws.Range("G26").Value = 40
$ws.Range("H26").Value = 35
$ws.Range("I26").Value = 518
$ws.Range("J26").Value = 475
$ws.Range("K26").Value = 9.052631578947
$ws.Range("L26").Value = 29.5
$ws.Range("M26").Value = -9.122807017543

# Row 27 - UCR Rape* (D27 becomes "0" placeholder, E27 becomes "***.*")
$ws.Range("D27").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 700
$ws.Range("I27").Value = 32
$ws.Range("K27").Value = 39.130434782608
$ws.Range("L27").Value = -5.882352941176

# Row 28 - Other Sex Crimes
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = -44.444444444444
$ws.Range("I28").Value = 53
$ws.Range("J28").Value = 78
$ws.Range("K28").Value = -32.051282051282
$ws.Range("L28").Value = -14.516129032258

# Row 33 - Traffic Fatalities (C33 becomes "0" placeholder)
$ws.Range("C33").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C33").PasteSpecial(-4122)
